$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.994.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.238.43'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.44'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.23%  '
$ws.Range('E8').Value = '  +0.23%  '
$ws.Range('E9').Value = '  -3.29%  '
$ws.Range('E10').Value = '  -5.09%  '
$ws.Range('E11').Value = '  -1.86%  '
$ws.Range('E12').Value = '  -3.78%  '
$ws.Range('E13').Value = '  -0.70%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.579.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.292.58'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('E16').Value = '  -2.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.48'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.878.07'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0963'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.16'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.52%  '
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.51'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('E24').Value = '  -5.24%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.86%  '
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.20'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('B28').Value = 'InjectiveProtocol'
$ws.Range('C28').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.50%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.84'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.29%  '
$ws.Range('E30').Value = '  -1.14%  '
$ws.Range('E31').Value = '  -0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.67%  '
$ws.Range('E33').Value = '  -5.22%  '
$ws.Range('E34').Value = '  -2.69%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.17'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.120'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.76'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.21%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '14.97'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.46%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.04%  '
$ws.Range('E41').Value = '  -6.53%  '
$ws.Range('E42').Value = '  -5.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.738.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '84.73'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.30%  '
$ws.Range('E46').Value = '  -4.61%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '99.80'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.88%  '
$ws.Range('E48').Value = '  -3.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '68.78'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.12%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '14.25'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.92%  '
